# "added the error message on Login page"
#
# Adds a new "InvalidLogin" worksheet right after the existing
# "ValidLogin" sheet. It mirrors the UserName/Password header row and
# the "admin" username row, but uses "Atharv" as the (incorrect)
# password so the scenario demonstrates a failed login.

$wb = $excel.ActiveWorkbook
$validLogin = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after "ValidLogin".
$invalidLogin = $wb.Worksheets.Add($null, $validLogin, 1)
$invalidLogin.Name = "InvalidLogin"

$invalidLogin.Range("A1").Value = "UserName"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("A2").Value = "admin"
$invalidLogin.Range("B2").Value = "Atharv"

# Make the new sheet the active tab/selection, like it was left after
# the edit.
$invalidLogin.Activate()
$invalidLogin.Range("B4").Select()
